$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9301554646195587
$ws.Range("C2").Value = 0.6058669934699258
$ws.Range("D2").Value = 0.7337789780881485
$ws.Range("B3").Value = 0.8686266612010792
$ws.Range("C3").Value = 0.9465721429073342
$ws.Range("D3").Value = 0.9059258944183849
$ws.Range("B4").Value = 0.5800353423767848
$ws.Range("C4").Value = 0.5509822308055513
$ws.Range("D4").Value = 0.5651356344350986
$ws.Range("B5").Value = 0.6801786589471244
$ws.Range("C5").Value = 0.9011592962120328
$ws.Range("D5").Value = 0.7752287480300734
$ws.Range("B7").Value = 0.06289308176100629
$ws.Range("C7").Value = 0.0001186394428691763
$ws.Range("D7").Value = 0.0002368321333838575
$ws.Range("B8").Value = 0.4716636197440585
$ws.Range("C8").Value = 0.001205584942337526
$ws.Range("D8").Value = 0.00240502258204343
$ws.Range("B11").Value = 0.7865976375601036
$ws.Range("C11").Value = 0.7297303430583276
$ws.Range("D11").Value = 0.7570976348685434
$ws.Range("B12").Value = 0.2425787019606707
$ws.Range("C12").Value = 0.04159197556918341
$ws.Range("D12").Value = 0.07100892698186628
$ws.Range("B13").Value = 0.8792954453401336
$ws.Range("C13").Value = 0.8940070379505203
$ws.Range("D13").Value = 0.8865902168175908
$ws.Range("B14").Value = 0.3403777526329533
$ws.Range("C14").Value = 0.1655490644314622
$ws.Range("D14").Value = 0.2227564011278517
$ws.Range("B16").Value = 0.6240941177056816
$ws.Range("C16").Value = 0.8412528493885622
$ws.Range("D16").Value = 0.7165824430614578
$ws.Range("B22").Value = 0.7739371222418707
$ws.Range("C22").Value = 0.7739371222418707
$ws.Range("D22").Value = 0.7739371222418707
$ws.Range("E22").Value = 0.7739371222418707
$ws.Range("B23").Value = 0.3233248241924577
$ws.Range("C23").Value = 0.2839018079089053
$ws.Range("D23").Value = 0.2818373366272222
$ws.Range("B24").Value = 0.7506928016837868
$ws.Range("C24").Value = 0.7739371222418707
$ws.Range("D24").Value = 0.7493575289337119
